$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scores (column E) for several students
$ws.Range("E2").Value = 86
$ws.Range("E3").Value = 78
$ws.Range("E5").Value = 81.3
$ws.Range("E6").Value = 72.3
$ws.Range("E7").Value = 71

# Add committee review comments (column I) for rows 2 and 5
$comment = "该生毕业设计选题具有实际工程意义，设计质量良好，设计成果完整，具有较好的工程价值。答辩准备充分、称述问题清晰、回答问题良好。"
$ws.Range("I2").Value = $comment
$ws.Range("I5").Value = $comment

# Update the selected cell in the sheet view
$ws.Range("I5").Select()
